# "Generate Report for Handoff"
#
# The localization-status report moved from "In Translation" to
# "Ready for handoff": the Status column on every sheet is refreshed,
# along with the handoff timestamps that get stamped when the report is
# (re)generated. The Status/Datetime columns also widen slightly to fit
# the new text.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
$ovw.Range("E2").Value = "Ready for handoff"   # zh-cn status (Overview)
$ovw.Range("F2").Value = "Ready for handoff"   # de-de status (Overview)
$zh.Range("C2").Value  = "Ready for handoff"   # Status column on zh-cn sheet
$de.Range("C2").Value  = "Ready for handoff"   # Status column on de-de sheet

# --- Handoff timestamps, refreshed for the new handoff generation ----
$ovw.Range("G2").Value = "2016-08-27 08:58:04"  # Latest HO Xliff Generate Date
$de.Range("H2").Value  = "2016-08-27 08:58:04"  # Latest Handoff Datetime (de-de)
$zh.Range("H2").Value  = "2016-08-27 08:57:57"  # Latest Handoff Datetime (zh-cn)

# --- Widen the Status / Datetime columns to fit the new text ---------
$ovw.Columns.Item(5).ColumnWidth = 16.33   # Overview!E (zh-cn status)
$ovw.Columns.Item(6).ColumnWidth = 16.33   # Overview!F (de-de status)
$zh.Columns.Item(3).ColumnWidth  = 16.33   # zh-cn!C (Status)
$de.Columns.Item(3).ColumnWidth  = 16.33   # de-de!C (Status)
